# Update job_info.xlsx: replace the RDMSTG/CMNREF job rows with the new
# D800 Customer/BMSIW job rows, and add the two hyperlinks that go with them.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Wipe out all the old data rows (2-7) below the header row; the new
# content only needs rows 2-4.
$ws.Range("A2:F7").Clear()

# Row 3 - Parallel job with its job/control ids (set before the hyperlinked
# A3 text so the shared-string table is built in the same order as the
# original authoring session).
$ws.Range("B3").Value = "Parallel"
$ws.Range("C3").Value = "J800801"
$ws.Range("D3").Value = "C800801"

# Row 2 - sequence job, plus the dependency reference on row 3.
$ws.Range("A2").Value = "LD_D800_Customer_JobSeq"
$ws.Range("B2").Value = "Sequence"
$ws.Range("F3").Value = "LD_D800_Customer_JobSeq"

# A3 carries a hyperlink to the customer job reference page; set the link
# first (it seeds the cell with the display text) then overwrite the cell
# text with the real job name while keeping the hyperlink + display text.
$ws.Hyperlinks.Add($ws.Range("A3"), "http://iwrefresh.w3ibm.mybluemix.net/Domains/ODS ADL/Jobs/LD_RDHIW_CUSTOMER_From_BMSIW_PJob@", "", "", "LD_RDHIW_CUSTOMER_From_BMSIW_PJob@")
$ws.Range("A3").Value = "LD_D800_JAPAN_CUSTOMER_REFERENCE_from_BMSIW_Delta_NZ_PJob"

# E2 is the datagroup hyperlink/label.
$ws.Hyperlinks.Add($ws.Range("E2"), "http://iwrefresh.w3ibm.mybluemix.net/Domains/ODS ADL/Datagroups/D800 CUSTOMER-BMSIW ETL 01 D", "", "", "D800 CUSTOMER-BMSIW ETL 01 D")
$ws.Range("E2").Value = "D800 CUSTOMER-BMSIW ETL 01 D"

# E4 stays empty but keeps the wrap-text style used elsewhere in the sheet.
$ws.Range("E4").WrapText = $true

# Restore the selection to F7 as in the final saved state.
[void]$ws.Range("F7").Select()
